$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processes")

# The 'type' column (PRODUCTIVE/DISSIPATIVE, used for Process Type data validation)
# and the 'description' column were swapped: 'type' moved from column E to column B,
# and 'description' moved from column B to column E.
for ($r = 1; $r -le 7; $r++) {
    $bVal = $ws.Cells.Item($r, 2).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 2).Value = $eVal
    $ws.Cells.Item($r, 5).Value = $bVal
}

# Column widths follow the swapped content (best fit).
$ws.Columns.Item(2).ColumnWidth = 13.833333
$ws.Columns.Item(5).ColumnWidth = 26.0

# The named range for the CGAM example table no longer includes the (now-moved)
# fifth column.
$wb.Names.Item("cgam_processes").RefersTo = "=Processes!`$A`$1:`$D`$4"

# Update the saved selection/active cell on the Processes sheet.
$ws.Activate()
$ws.Range("D4").Select() | Out-Null
